# Fruta / hortaliza, semanal
# This update re-shuffles the weekly price records (rows 2-21) by moving the
# Fecha/Volumen/Precio minimo/Precio maximo/Precio promedio/Origen/Precio $Kg
# values (columns D, J, K, L, M, O, P) between rows, keeping the rest of each
# row (market/region/category identity columns) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (where the "new" values currently live)
$map = @{
    2  = 6
    3  = 19
    4  = 17
    5  = 12
    6  = 3
    7  = 9
    8  = 14
    9  = 7
    10 = 8
    11 = 11
    12 = 21
    13 = 16
    14 = 15
    15 = 2
    16 = 18
    17 = 4
    18 = 20
    19 = 10
    20 = 5
    21 = 13
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the relevant columns before overwriting
# anything, since sources and destinations overlap.
# NOTE: use Value2 for reading (Value's getter is unreliable in this
# runtime), Value works fine for writing.
$snapshot = @{}
for ($r = 2; $r -le 21; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
